# Slide 12 ("Model Selection (Demand)") - Content Placeholder 2
#
# Hunk 1: collapse the 3 runs "Random Forest " / "&gt; Linear " /
#         "Regression (once again)" into a single run reading
#         "Random Forest > Linear Regression (once again)".
#
# Hunk 2: split the run "Linear Regression: 0.61" into two runs,
#         "Linear Regression" and ": 0.71" (score bumped 0.61 -> 0.71).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Hunk 1 ---------------------------------------------------------
$para1 = $tr.Paragraphs(1)

# Runs 3, 4, 5 of this paragraph are:
#   3: "Random Forest "
#   4: "&gt; Linear "
#   5: "Regression (once again)"
# Remove the trailing two runs (highest index first, so earlier
# indices stay stable), keeping run 3's formatting, then give run 3
# the fully merged text.
$para1.Runs(5).Text = ""
$para1.Runs(4).Text = ""
$para1.Runs(3).Text = "Random Forest > Linear Regression (once again)"

# --- Hunk 2 ---------------------------------------------------------
$para3 = $tr.Paragraphs(3)
$run1  = $para3.Runs(1)

# Replace the trailing ": 0.61" (characters 18-23 of the paragraph)
# with ": 0.71", then trim run 1's text back down to "Linear
# Regression" so the paragraph ends up as two runs.
$tail = $run1.Characters(18, 6)
$tail.Text = ": 0.71"
$run1.Text = "Linear Regression"
